$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Renumber the "id" column (process ids) from the old 1-based / offset
# numbering to a fresh 0-based sequential numbering.
$ws.Range("S2").Value = 0
$ws.Range("S3").Value = 1
$ws.Range("S4").Value = 2
$ws.Range("S5").Value = 3
$ws.Range("S6").Value = 4
$ws.Range("S7").Value = 5

# Reflect the cell range the author was working on (the id column) as the
# active selection on the params sheet.
$ws.Activate()
$ws.Range("S2:S7").Select()
